# Automatische test-sync: 2025-07-27 18:32:50
# Adds the second test-mail row to the "Logs" sheet and refreshes the
# "Aantal" counter for "Overig" on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append row 3 --------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Zorg jij dat dit geregeld wordt?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #2: Zorg jij dat dit geregeld wordt?"
$logs.Range("D3").Value = "Overig"
$logs.Range("E3").Value = "Beste klant,`nBedankt voor uw e-mail. Om u beter te kunnen helpen, zou u meer specifieke informatie kunnen geven over welk aspect van onze service u graag geregeld zou willen hebben? Zo kunnen wij u adequaat assisteren. Alvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Range("F3").Value = "2025-07-27 18:32:43"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Ja"
$logs.Range("J3").Value = "Nee"

# The multi-line text in E3 makes Excel auto-wrap & auto-grow row 3; re-run
# AutoFit so the row goes back to the (unpinned) standard height, matching
# a row that never had an explicit height set.
$logs.Rows.Item(3).AutoFit()

# ---- Logs sheet: stretch the existing conditional formatting down to the
#      newly added row 3 (one ModifyAppliesToRange call per rule-group, since
#      all cfRules sharing a sqref move together) --------------------------
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))
$logs.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J3"))

# ---- Dashboard sheet: bump the "Overig" count from 1 to 2 ----------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 2
